$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.932.51'
$ws.Range('E2').Value = '  -5.46%  '
$ws.Range('D3').Value = '3.017.75'
$ws.Range('E3').Value = '  -7.81%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '546.50'
$ws.Range('E5').Value = '  -7.92%  '
$ws.Range('D6').Value = '136.55'
$ws.Range('E6').Value = '  -9.73%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.008.52'
$ws.Range('E8').Value = '  -7.86%  '
$ws.Range('D9').Value = '0.481'
$ws.Range('E9').Value = '  -11.85%  '
$ws.Range('D10').Value = '6.27'
$ws.Range('E10').Value = '  -6.96%  '
$ws.Range('D11').Value = '0.150'
$ws.Range('E11').Value = '  -12.85%  '
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').Value = '  -11.02%  '
$ws.Range('D13').Value = '35.23'
$ws.Range('E13').Value = '  -8.66%  '
$ws.Range('D14').Value = '0.0000216'
$ws.Range('E14').Value = '  -12.85%  '
$ws.Range('D15').Value = '3.500.89'
$ws.Range('E15').Value = '  -7.76%  '
$ws.Range('D16').Value = '64.003.27'
$ws.Range('E16').Value = '  -5.30%  '
$ws.Range('D17').Value = '0.110'
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('D18').Value = '3.021.48'
$ws.Range('E18').Value = '  -7.57%  '
$ws.Range('D19').Value = '6.47'
$ws.Range('E19').Value = '  -9.48%  '
$ws.Range('D20').Value = '474.60'
$ws.Range('E20').Value = '  -11.21%  '
$ws.Range('D21').Value = '13.39'
$ws.Range('E21').Value = '  -10.82%  '
$ws.Range('D22').Value = '0.672'
$ws.Range('E22').Value = '  -11.54%  '
$ws.Range('D23').Value = '6.96'
$ws.Range('E23').Value = '  -12.01%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '77.29'
$ws.Range('E24').Value = '  -9.76%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '12.26'
$ws.Range('E25').Value = '  -9.82%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '2.67'
$ws.Range('E27').Value = '  -17.91%  '
$ws.Range('D28').Value = '2.05'
$ws.Range('E28').Value = '  -5.50%  '
$ws.Range('D29').Value = '7.60'
$ws.Range('E29').Value = '  -6.42%  '
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '579.17'
$ws.Range('E31').Value = '  +11.53%  '
$ws.Range('D32').Value = '2.60'
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '25.38'
$ws.Range('E33').Value = '  -13.33%  '
$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').Value = '1.10'
$ws.Range('E34').Value = '  -6.20%  '
$ws.Range('D35').Value = '5.27'
$ws.Range('E35').Value = '  -8.52%  '
$ws.Range('D36').Value = '5.74'
$ws.Range('E36').Value = '  -13.84%  '
$ws.Range('D37').Value = '51.40'
$ws.Range('E37').Value = '  -3.75%  '
$ws.Range('D38').Value = '0.0411'
$ws.Range('E38').Value = '  -7.41%  '
$ws.Range('D39').Value = '0.0781'
$ws.Range('E39').Value = '  -9.10%  '
$ws.Range('D40').Value = '0.118'
$ws.Range('E40').Value = '  -6.22%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.69'
$ws.Range('E41').Value = '  -3.90%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.902.03'
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '8.09'
$ws.Range('E43').Value = '  -10.07%  '
$ws.Range('D45').Value = '0.237'
$ws.Range('E45').Value = '  -11.45%  '
$ws.Range('D46').Value = '2.03'
$ws.Range('E46').Value = '  -7.89%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '24.41'
$ws.Range('E47').Value = '  -9.01%  '
$ws.Range('B48').Value = 'PEPE'
$ws.Range('C48').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D48').Value = '0.0₃0520'
$ws.Range('E48').Value = '  -12.14%  '
$ws.Range('D49').Value = '117.40'
$ws.Range('E49').Value = '  -5.22%  '
$ws.Range('D50').Value = '0.107'
$ws.Range('E50').Value = '  -6.56%  '
$ws.Range('D51').Value = '2.00'
$ws.Range('E51').Value = '  -14.51%  '
